$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary area updates ---------------------------------------
$ws.Range("E11").Value = 227760          # VALOR MORA total (56940 * 4 rows)
$ws.Range("C13").Value = 2               # Cant. Trabajadores
$ws.Range("F13").Value = 3               # Cant. Periodos

# --- Insert two additional data rows below the existing table ------------
# Before: row16 = GUSTAVO... (now), row17 = JORGE... (was the last/bottom-border row)
# After : rows 16-18 use the "middle" style, row19 becomes the new bottom row.
$ws.Rows("18:19").Insert()

# Row 19 should inherit the bold bottom-border ("last row") formatting that
# used to belong to row 17.
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# Rows 17 and 18 should look like the existing "middle" row (row16).
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$ws.Range("B18:J18").PasteSpecial(-4122)

# --- Populate the new / shifted data rows ---------------------------------
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73214230"
$ws.Range("D17").Value = "JORGE ENRIQUE PEÑA GUZMAN"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73214230"
$ws.Range("D18").Value = "JORGE ENRIQUE PEÑA GUZMAN"
$ws.Range("E18").Value = "2506"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73214230"
$ws.Range("D19").Value = "JORGE ENRIQUE PEÑA GUZMAN"
$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# --- Cosmetic: widen column D slightly, matching the wider new name -------
$ws.Columns("D").ColumnWidth = 31.453125
